$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell $ws "D2" "62.798.57"
$ws.Range("E2").Value = "  -0.82%  "
Set-TextCell $ws "D3" "2.537.99"
$ws.Range("E3").Value = "  +3.28%  "
Set-TextCell $ws "D4" "1.00"
$ws.Range("E4").Value = "  +0.00%  "
Set-TextCell $ws "D5" "565.87"
$ws.Range("E5").Value = "  -0.52%  "
Set-TextCell $ws "D6" "146.76"
$ws.Range("E6").Value = "  +2.83%  "
Set-TextCell $ws "D7" "1.00"
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -1.95%  "
Set-TextCell $ws "D9" "2.537.88"
$ws.Range("E9").Value = "  +3.28%  "
Set-TextCell $ws "D10" "0.104"
$ws.Range("E10").Value = "  -1.32%  "
$ws.Range("E11").Value = "  -2.94%  "
$ws.Range("E12").Value = "  +0.51%  "
Set-TextCell $ws "D13" "0.351"
$ws.Range("E13").Value = "  -0.72%  "
Set-TextCell $ws "D14" "26.88"
$ws.Range("E14").Value = "  +1.98%  "
Set-TextCell $ws "D15" "2.997.96"
$ws.Range("E15").Value = "  +3.43%  "
Set-TextCell $ws "D16" "62.800.85"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("E17").Value = "  -1.82%  "
Set-TextCell $ws "D18" "2.535.13"
$ws.Range("E18").Value = "  +2.93%  "
Set-TextCell $ws "D19" "11.44"
$ws.Range("E19").Value = "  +1.42%  "
Set-TextCell $ws "D20" "333.17"
$ws.Range("E20").Value = "  -2.69%  "
$ws.Range("E21").Value = "  -1.43%  "
Set-TextCell $ws "D22" "6.74"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("E23").Value = "  +0.04%  "
Set-TextCell $ws "D24" "64.68"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("E25").Value = "  -3.39%  "
Set-TextCell $ws "D26" "1.58"
$ws.Range("E26").Value = "  +3.40%  "
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  +10.87%  "
Set-TextCell $ws "D29" "8.29"
$ws.Range("E29").Value = "  +0.99%  "
Set-TextCell $ws "D30" "7.22"
$ws.Range("E30").Value = "  +5.09%  "
Set-TextCell $ws "D31" "0.0₃0806"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("E32").Value = "  -0.54%  "
Set-TextCell $ws "D33" "176.63"
$ws.Range("E33").Value = "  +0.84%  "
Set-TextCell $ws "D34" "1.57"
$ws.Range("E34").Value = "  +4.01%  "
Set-TextCell $ws "D35" "404.07"
$ws.Range("E35").Value = "  +9.20%  "
Set-TextCell $ws "D36" "0.395"
$ws.Range("E36").Value = "  -1.37%  "
Set-TextCell $ws "D37" "18.76"
$ws.Range("E37").Value = "  -1.03%  "
Set-TextCell $ws "D39" "4.30"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("E40").Value = "  +0.40%  "
Set-TextCell $ws "D41" "1.00"
$ws.Range("E41").Value = "  +0.06%  "
Set-TextCell $ws "D42" "38.96"
$ws.Range("E42").Value = "  -3.47%  "
Set-TextCell $ws "D43" "151.40"
$ws.Range("E43").Value = "  +0.89%  "
Set-TextCell $ws "D44" "3.72"
$ws.Range("E44").Value = "  -0.19%  "
Set-TextCell $ws "D45" "20.49"
$ws.Range("E45").Value = "  -1.13%  "
Set-TextCell $ws "D46" "0.601"
$ws.Range("E46").Value = "  +0.35%  "
Set-TextCell $ws "D47" "0.0955"
$ws.Range("E47").Value = "  -1.08%  "
Set-TextCell $ws "D48" "0.0516"
$ws.Range("E48").Value = "  -2.00%  "
Set-TextCell $ws "D49" "0.0234"
$ws.Range("E49").Value = "  +3.45%  "
Set-TextCell $ws "D50" "18.23"
$ws.Range("E50").Value = "  +0.42%  "
Set-TextCell $ws "D51" "1.75"
$ws.Range("E51").Value = "  +0.71%  "
